$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Typo fix: "Rende sample error table" -> "Render sample error table"
#    (hidden inside a paragraph about removing validateSamp code)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Rende sample error table", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Render sample error table", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Append four new bullet paragraphs at the end of the document describing
#    the 1/24/2024 bug-fix update log entry.
# ---------------------------------------------------------------------------

# -- Paragraph 1: top level bullet (same level as the other dated entries) --
$endRange = $d.Paragraphs.Last.Range
$endRange.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$r1 = $p1.Range
$r1.Text = "1/24/2024 – Dan bug fixes found by Kendal Robbins related to new “verified” columns not working"
$r1.ListFormat.ListLevelNumber = 1

# -- Paragraph 2: sub-bullet --
$endRange = $d.Paragraphs.Last.Range
$endRange.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$r2 = $p2.Range
$r2.Text = "Had some “or” statements where I needed “and”"
$r2.ListFormat.ListLevelNumber = 2

# -- Paragraph 3: sub-bullet --
$endRange = $d.Paragraphs.Last.Range
$endRange.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$r3 = $p3.Range
$r3.Text = "Did not anticipate people writing “verified” with quote marks (but I can see from the instructions why they would...I also changed instructions to indicate quotes not needed), so added code to pick up “verified” with quotes as marked verified.  "
$r3.ListFormat.ListLevelNumber = 2

# -- Paragraph 4: sub-bullet --
$endRange = $d.Paragraphs.Last.Range
$endRange.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$r4 = $p4.Range
$r4.Text = "Found/fixed code issue in the non_integerTL code I had...was minor but produced an error if you clicked on the button to get row numbers when there were no offending rows."
$r4.ListFormat.ListLevelNumber = 2

Write-Host "Edit complete"
